# ---------------------------------------------------------------------------
# BenthFun Spring 2023 diving log - apply commit:
#  "Anto Ali ALK LOW + Exp Amb Jerem + 40% ALK AMB"
#
#  1) Insert a new empty worksheet named "Sheet2" right before "Sheet1".
#  2) Append 14 new observation rows (195-208) to the "Corrected" sheet
#     for the new "AMB" (ambient) incubation chambers.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1) Insert the new, empty "Sheet2" just before the existing "Sheet1" ---
$sheet1 = $wb.Worksheets.Item("Sheet1")
$newSheet = $wb.Worksheets.Add($sheet1)
# Add() names it "Sheet2" automatically (next free default name); nothing else to do.

# --- 2) Add the new rows to the "Corrected" worksheet ----------------------
$ws = $wb.Worksheets.Item("Corrected")
$ws.Activate()

$rows = @(
    @{ r=195; A=45092; B="Tn_t1_AMB_tile_01"; C="Light_01"; D=0.4381944444444445; E=0.4847222222222222; F=0.4847222222222222; G=1; H="A1"; I=24.99; J=-74; Bstyle="plain" },
    @{ r=196; A=45092; B="Tn_t1_AMB_tile_02"; C="Light_02"; D=0.44027777777777777; E=0.48819444444444443; F=0.48819444444444443; G=2; H="A2"; I=24.67; J=-74.1; Bstyle="plain" },
    @{ r=197; A=45092; B="Tn_t1_AMB_tile_03"; C="Light_03"; D=0.44097222222222227; E=0.4909722222222222; F=0.4909722222222222; G=3; H="A3"; I=24.74; J=-71.4; Bstyle="plain" },
    @{ r=198; A=45092; B="Tn_t1_AMB_blank_01"; C="Light_04"; D=0.44236111111111115; E=0.49444444444444446; F=0.49444444444444446; G=4; H=$null; I=24.43; J=-68.8; Bstyle="plain" },
    @{ r=199; A=45092; B=$null; C="Dark_01"; D=0.4479166666666667; E=0.48055555555555557; F=0.48055555555555557; G=5; H="A7"; I=24; J=-64.6; Bstyle=$null },
    @{ r=200; A=45092; B=$null; C="Dark_02"; D=0.44930555555555557; E=0.48125; F=0.48125; G=6; H="A8"; I=23.99; J=-64.4; Bstyle=$null },
    @{ r=201; A=45092; B=$null; C="Dark_03"; D=0.45069444444444445; E=0.48194444444444445; F=0.48194444444444445; G=7; H="A9"; I=24; J=-62.2; Bstyle=$null },
    @{ r=202; A=45092; B="Tn_t1_AMB_tile_04"; C="Light_05"; D=0.5048611111111111; E=0.548611111111111; F=0.548611111111111; G=5; H="A7"; I=25.19; J=-72.2; Bstyle="font2" },
    @{ r=203; A=45092; B="Tn_t1_AMB_tile_05"; C="Light_06"; D=0.5055555555555555; E=0.5520833333333334; F=0.5520833333333334; G=6; H="A8"; I=24.99; J=-69.9; Bstyle="font2" },
    @{ r=204; A=45092; B="Tn_t1_AMB_tile_06"; C="Light_07"; D=0.50625; E=0.5555555555555556; F=0.5555555555555556; G=7; H="A9"; I=25.17; J=-63.6; Bstyle="font2" },
    @{ r=205; A=45092; B="Tn_t1_AMB_blank_02"; C="Light_08"; D=0.5069444444444444; E=0.5583333333333333; F=0.5583333333333333; G=4; H=$null; I=24.85; J=-67.1; Bstyle="plain" },
    @{ r=206; A=45092; B=$null; C="Dark_05"; D=0.5125000000000001; E=0.5437500000000001; F=0.5437500000000001; G=1; H="A1"; I=24.45; J=-63.3; Bstyle=$null },
    @{ r=207; A=45092; B=$null; C="Dark_06"; D=0.5131944444444444; E=0.545138888888889; F=0.545138888888889; G=2; H="A2"; I=24.33; J=-62.1; Bstyle=$null },
    @{ r=208; A=45092; B=$null; C="Dark_07"; D=0.513888888888889; E=0.5458333333333333; F=0.5458333333333333; G=3; H="A3"; I=24.37; J=-59.5; Bstyle=$null }
)

# Reference cells whose formats we reuse (none of these are themselves edited).
$dateFmtSrc = $ws.Range("A194")   # numFmtId 14 source -> we instead set a fresh "d-mmm" format below
$timeFmtSrc = $ws.Range("D194")   # custom time format (style index 7)
$numFmtSrc  = $ws.Range("I194")   # 0.00 numeric format (style index 3)
$rightSrc   = $ws.Range("H194")   # right-aligned tile-number style (style index 9)
$fontSrc    = $ws.Range("B10")    # alternate font style used for some corrected labels (style index 2)

# The workbook's shared-string table is appended to in the order cells are
# *filled*. The source edit filled the sheet column-by-column (Label down the
# whole range, then Tile_N степ down the whole range, ...), so we replay the
# same column-major order here to land identical shared-string ids.

# --- Column A: Diving_Date ---
foreach ($row in $rows) {
    $r = $row.r
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 1).NumberFormat = "d-mmm"
}

# --- Column B: Label (optional) ---
foreach ($row in $rows) {
    if ($row.B -ne $null) {
        $r = $row.r
        $ws.Cells.Item($r, 2).Value = $row.B
    }
}

# --- Column C: Chamber ---
foreach ($row in $rows) {
    $r = $row.r
    $ws.Cells.Item($r, 3).Value = $row.C
}

# --- Columns D/E/F: Start_incubation / Stop_Incubation / Stop_Alkalinity ---
foreach ($row in $rows) {
    $r = $row.r
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
}

# --- Column G: O2_sensor_used ---
foreach ($row in $rows) {
    $r = $row.r
    $ws.Cells.Item($r, 7).Value = $row.G
}

# --- Column H: Tile_N deg (optional) ---
foreach ($row in $rows) {
    if ($row.H -ne $null) {
        $r = $row.r
        $ws.Cells.Item($r, 8).Value = $row.H
    }
}

# --- Columns I/J: Temperature / pH_mV ---
foreach ($row in $rows) {
    $r = $row.r
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
}

# --- Now apply number formats / font overrides, again per affected range ---
foreach ($row in $rows) {
    $r = $row.r

    if ($row.B -ne $null -and $row.Bstyle -eq "font2") {
        $fontSrc.Copy()
        $ws.Cells.Item($r, 2).PasteSpecial(-4122)
    }

    $timeFmtSrc.Copy()
    $ws.Range($ws.Cells.Item($r, 4), $ws.Cells.Item($r, 6)).PasteSpecial(-4122)

    if ($row.H -ne $null) {
        $rightSrc.Copy()
        $ws.Cells.Item($r, 8).PasteSpecial(-4122)
    }

    $numFmtSrc.Copy()
    $ws.Range($ws.Cells.Item($r, 9), $ws.Cells.Item($r, 10)).PasteSpecial(-4122)
}

# Move the active selection the way the source workbook shows it ("J209" is the
# next free cell below the appended data).
$ws.Range("J209").Select()

Write-Output "done"
